$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.976.44"
$ws.Range("E2").Value = "  +4.32%  "
$ws.Range("D3").Value = "3.077.92"
$ws.Range("E3").Value = "  +2.59%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'580.00"
$ws.Range("E5").Value = "  +3.07%  "
$ws.Range("D6").Value = "'142.22"
$ws.Range("E6").Value = "  +2.89%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.069.31"
$ws.Range("E8").Value = "  +2.73%  "
$ws.Range("D9").Value = "'0.530"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("E10").Value = "  +5.72%  "
$ws.Range("D11").Value = "'5.71"
$ws.Range("E11").Value = "  +11.54%  "
$ws.Range("E12").Value = "  +2.26%  "
$ws.Range("E13").Value = "  +4.47%  "
$ws.Range("E14").Value = "  +4.55%  "
$ws.Range("D15").Value = "'0.124"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").Value = "3.586.83"
$ws.Range("E16").Value = "  +2.86%  "
$ws.Range("D17").Value = "'7.28"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "3.075.70"
$ws.Range("E18").Value = "  +2.83%  "
$ws.Range("D19").Value = "61.895.24"
$ws.Range("E19").Value = "  +4.49%  "
$ws.Range("D20").Value = "'448.16"
$ws.Range("E20").Value = "  +4.20%  "
$ws.Range("D21").Value = "'14.02"
$ws.Range("E21").Value = "  +2.69%  "
$ws.Range("D22").Value = "'0.735"
$ws.Range("E22").Value = "  +2.07%  "
$ws.Range("E23").Value = "  +4.65%  "
$ws.Range("D24").Value = "'13.75"
$ws.Range("E24").Value = "  +3.21%  "
$ws.Range("D25").Value = "'81.67"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("E27").Value = "  +4.56%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("E29").Value = "  +4.54%  "
$ws.Range("D30").Value = "'8.22"
$ws.Range("E30").Value = "  +5.52%  "
$ws.Range("D31").Value = "'6.83"
$ws.Range("E31").Value = "  +11.61%  "
$ws.Range("D32").Value = "'0.113"
$ws.Range("E32").Value = "  +14.32%  "
$ws.Range("D33").Value = "'26.86"
$ws.Range("E33").Value = "  +4.41%  "
$ws.Range("E34").Value = "  +4.58%  "
$ws.Range("D35").Value = "0.0₃0796"
$ws.Range("E35").Value = "  +4.00%  "
$ws.Range("D36").Value = "'6.06"
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("D37").Value = "'2.21"
$ws.Range("E37").Value = "  +5.75%  "
$ws.Range("D38").Value = "'50.18"
$ws.Range("E38").Value = "  +2.45%  "
$ws.Range("E39").Value = "  +8.63%  "
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("D41").Value = "'418.91"
$ws.Range("E41").Value = "  +4.47%  "
$ws.Range("D42").Value = "2.924.99"
$ws.Range("E42").Value = "  +6.01%  "
$ws.Range("E43").Value = "  +5.14%  "
$ws.Range("E44").Value = "  +9.78%  "
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("E46").Value = "  +6.51%  "
$ws.Range("D48").Value = "'34.88"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").Value = "'123.74"
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").Value = "'24.45"
